$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header labels
$ws.Range("A1").Value = "button_testResultActions_class"
$ws.Range("B1").Value = "button_testResultActions_class_1"
$ws.Range("C1").Value = "button_testResultActions_internalRoleButtonName"
$ws.Range("D1").Value = "button_testResultActions_internalRoleButtonName_1"
$ws.Range("E1").Value = "button_testResultDetails_class"
$ws.Range("F1").Value = "button_testResultDetails_internalRoleButtonName"

# Update row 2 data values (C2, D2, E2 shift/change; A2, B2, F2 stay the same)
$ws.Range("C2").Value = "Failed Automations - Apply to"
$ws.Range("D2").Value = "Failed Portal - Login with"
$ws.Range("E2").Value = '"]:nth-child(3) [class="css-1yjo05o'

# Adjust column widths to match new layout: A=37 B=37 C=49 D=51 E=37 F=80
# (Excel's ColumnWidth COM property is offset from the raw OOXML stored width
# by 5/6 of a character, so subtract that padding to land on the exact value.)
$pad = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 37 - $pad
$ws.Columns.Item(2).ColumnWidth = 37 - $pad
$ws.Columns.Item(3).ColumnWidth = 49 - $pad
$ws.Columns.Item(4).ColumnWidth = 51 - $pad
$ws.Columns.Item(5).ColumnWidth = 37 - $pad
$ws.Columns.Item(6).ColumnWidth = 80 - $pad
